# Generate Report for Handback
#
# The file "d56253df-ed40-41db-8663-8ed10e624083.md" has now been handed
# back successfully (it was previously "Ready for handoff"/pending), so
# update the status + handback metadata across all three report sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
# Row 3 is the d56253df... file; its zh-cn / de-de status columns move
# from "Ready for handoff" to "Handed back: in sync with en-US".
$ovw.Range("E3").Value = "Handed back: in sync with en-US"
$ovw.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
# Row 3 is the d56253df... file.
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-12 04:50:58"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).ColumnWidth = 12.83

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
# Row 3 is the d56253df... file.
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-12 04:51:09"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).ColumnWidth = 12.83
